# Applies the "Generate Report for Handback" edit:
# - refreshes handback UUIDs/hash/timestamps for two source files across
#   the Overview, zh-cn and de-de sheets, and updates the matching
#   hyperlink display text (while preserving the existing hyperlink targets).
$wb = $excel.ActiveWorkbook

# --- Overview sheet: update cell values ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws1.Range("B2").Value = "e2e\93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws1.Range("G2").Value = "2016-08-20 15:04:05"
$ws1.Range("A3").Value = "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws1.Range("B3").Value = "e2e\ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws1.Range("G3").Value = "2016-08-20 15:04:05"

# --- zh-cn sheet: update cell values ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws2.Range("G2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-20 15:03:58"
$ws2.Range("I2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws2.Range("J2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-20 15:04:26"
$ws2.Range("A3").Value = "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws2.Range("G3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-20 15:03:58"
$ws2.Range("I3").Value = "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws2.Range("J3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-20 15:04:26"

# --- de-de sheet: update cell values ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws3.Range("G2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-20 15:04:05"
$ws3.Range("I2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$ws3.Range("J2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-20 15:04:32"
$ws3.Range("A3").Value = "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws3.Range("G3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-20 15:04:05"
$ws3.Range("I3").Value = "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md"
$ws3.Range("J3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-20 15:04:32"

# --- Refresh hyperlink display text on all three sheets. ---
# The underlying link targets are unchanged by this edit, so the existing
# hyperlinks are removed and re-added in the same order (pointing at the
# same external URLs) purely to refresh the visible "display" text.

# Overview
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "e2e\93808e0a-e246-4825-aff9-e47cfeec904e.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "e2e\ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md")

# zh-cn
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "93808e0a-e246-4825-aff9-e47cfeec904e.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/633dab0429f84e4ab8ea583608553d2629847442/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "93808e0a-e246-4825-aff9-e47cfeec904e.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/633dab0429f84e4ab8ea583608553d2629847442/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md")

# de-de
$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "93808e0a-e246-4825-aff9-e47cfeec904e.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9d14160292c4ef818b14236fe698a9a86116aeeb/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "93808e0a-e246-4825-aff9-e47cfeec904e.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9d14160292c4ef818b14236fe698a9a86116aeeb/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "ffffc1163a5f-c534-41d2-a528-1889bbb7a4e9.md")
